# PhieuDeNghiTamUngTemplate.xlsx - "update exhibition report"
#
# Summary of the target edit (per the OOXML diff):
#   - workbook.xml: calcPr gets calcMode="manual"
#   - sheet1.xml:
#       * the two empty rows above the old "Tong cong" row are removed and
#         the "Tong cong" label ends up on row 15 (directly under the data
#         table), with everything below it (rows 18-21) shifting up to
#         rows 16-19
#       * the dimension / selected cell follow the new, smaller layout
#       * column F is widened to 32 characters

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level: switch calculation to manual (adds calcMode="manual") ---
$excel.Calculation = -4135   # xlCalculationManual

# --- Worksheet-level row cleanup ---
# Original layout:
#   row15 (blank) / row16 (blank) / row17 (D17 = "Tong cong") / row18.. (rest)
# Target layout:
#   row15 (D15 = "Tong cong") / row16.. (rest, shifted up by two)
#
# Write the "Tong cong" label onto row 15 first (it already sits in the
# shared-string table), then delete the two now-redundant rows above the
# old label row so everything below ripples up into place.
$ws.Range("D15").Value = "Tổng cộng"
$ws.Rows("16:16").Delete()
$ws.Rows("16:16").Delete()

# --- Column F width (18.85546875 -> 32) ---
$ws.Columns("F:F").ColumnWidth = 31.1

# --- Selection moves to F16 ---
$ws.Range("F16").Select()
